# "4er gewinnte Doku 2.0" - extend the documentation after the
# "Spielprinzip" bullet with the new outline / open-question paragraphs.

$d = $word.ActiveDocument

# Locate the "Spielprinzip" list-item paragraph robustly (by text, rather
# than a hard-coded index) and insert the new content right after it, i.e.
# right before the first of the two trailing blank paragraphs.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Spielprinzip") {
        $targetIndex = $i
    }
}

$insertionPara = $d.Paragraphs.Item($targetIndex + 1)
$insertionRange = $insertionPara.Range
$insertionRange.Collapse(1)

# Raw OOXML for the new paragraphs, built from small readable pieces and
# joined into one fragment. InsertXML replaces the (empty) contents of the
# collapsed range's paragraph, so the fragment below ends with one extra
# blank <w:p/> that takes the place of that original blank paragraph - the
# document keeps exactly the same two trailing blank paragraphs it had
# before, just pushed further down.
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$blankP = "<w:p $wns/>"

$pieces = @(
    # "2 Spieler" - new bullet in the same list as "Spielprinzip"
    "<w:p $wns>" +
        "<w:pPr><w:pStyle w:val=`"Listenabsatz`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr>" +
        "<w:r><w:t>2 Spieler</w:t></w:r>" +
    "</w:p>",
    $blankP,
    $blankP,

    "<w:p $wns><w:r><w:t>Betriebssystem:</w:t></w:r></w:p>",

    "<w:p $wns><w:r><w:t>C</w:t></w:r><w:r><w:t>#</w:t></w:r></w:p>",

    "<w:p $wns>" +
        "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>KonsolenApp</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`"> oder GUI?</w:t></w:r>" +
    "</w:p>",

    "<w:p $wns><w:r><w:t>Andere Programmiersprachen?</w:t></w:r></w:p>",

    $blankP,
    $blankP,
    $blankP,

    "<w:p $wns>" +
        "<w:r><w:t xml:space=`"preserve`">Welche </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>VScode</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/><w:proofErr w:type=`"gramStart`"/><w:r><w:t>erweiterungen</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t>?:</w:t></w:r><w:proofErr w:type=`"gramEnd`"/>" +
    "</w:p>"
)

# 13 further blank paragraphs, plus the one that stands in for the
# original trailing blank paragraph being replaced by InsertXML.
for ($i = 0; $i -lt 14; $i++) {
    $pieces += $blankP
}

$newParasXml = [string]::Join("", $pieces)

$insertionRange.InsertXML($newParasXml) | Out-Null
